$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 16 data (rounded trend row appended by calculator fix)
$ws.Range("A16").Value = 42622.887291666666

$ws.Range("B16").Value = -24
$ws.Range("C16").Value = 64
$ws.Range("D16").Value = 33
$ws.Range("E16").Value = 64
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 17897
$ws.Range("H16").Value = 10510
$ws.Range("I16").Value = 527
$ws.Range("J16").Value = 108
$ws.Range("K16").Value = 56
$ws.Range("L16").Value = 10
$ws.Range("M16").Value = 3
$ws.Range("N16").Value = "Named"
